$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.886.02'
$ws.Range("D3").Value = '1.740.34'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.10%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5181'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.46%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2748'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06164'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("D10").Value = '1.742.83'
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07173'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("B12").Value = 'Solana'
$ws.Range("C12").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.96'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6437'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.74%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.44'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.07%  '
$ws.Range("D18").Value = '25.901.42'
$ws.Range("E18").Value = '  +0.36%  '
$ws.Range("E19").Value = '  +2.55%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000006784'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.67%  '
$ws.Range("D21").Value = '1.963.84'
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.274'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.670'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.247'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '138.84'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.512'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.13'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.760'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '105.98'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.939'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08296'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.663'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04593'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.648'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9878'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.05%  '
$ws.Range("E36").Value = '  +0.68%  '
$ws.Range("E37").Value = '  +1.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01612'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.926'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9998'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '97.78'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.65%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.3834'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7391'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.981'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1128'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.205'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05262'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '54.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.41%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.598'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3402'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.11%  '
